$wb = $excel.ActiveWorkbook

# --- Section_A sheet ---
$ws = $wb.Worksheets.Item("Section_A")
$ws.Range("B2").Value = 'CS269'
$ws.Range("D2").Value = 'MA265'
$ws.Range("E2").Value = 'Free'
$ws.Range("B3").Value = 'CS269 (Tutorial)'
$ws.Range("C3").Value = 'CS265'
$ws.Range("E3").Value = 'CS265'
$ws.Range("F3").Value = 'CS268'
$ws.Range("C5").Value = 'Free'
$ws.Range("D5").Value = 'CS269'
$ws.Range("E5").Value = 'Free'
$ws.Range("F5").Value = 'CS269'
$ws.Range("B6").Value = 'CS268'
$ws.Range("C6").Value = 'Free'
$ws.Range("D6").Value = 'Free'
$ws.Range("E6").Value = 'CS268'
$ws.Range("F6").Value = 'CS265'
$ws.Range("C7").Value = 'Free'
$ws.Range("E7").Value = 'MA265'

# --- Section_B sheet ---
$ws = $wb.Worksheets.Item("Section_B")
$ws.Range("B2").Value = 'CS269 (Tutorial)'
$ws.Range("C2").Value = 'Free'
$ws.Range("D2").Value = 'Free'
$ws.Range("E2").Value = 'CS269'
$ws.Range("F2").Value = 'MA265'
$ws.Range("B3").Value = 'CS265'
$ws.Range("D3").Value = 'CS268'
$ws.Range("E3").Value = 'MA265'
$ws.Range("B5").Value = 'Free'
$ws.Range("C5").Value = 'CS269'
$ws.Range("D5").Value = 'CS265'
$ws.Range("F5").Value = 'CS265'
$ws.Range("D6").Value = 'CS269'
$ws.Range("F6").Value = 'CS268'
$ws.Range("B7").Value = 'Free'
$ws.Range("C7").Value = 'Free'
$ws.Range("E7").Value = 'CS268'

# --- Course_Summary sheet ---
$ws = $wb.Worksheets.Item("Course_Summary")
$ws.Range("A2").Value = 'MA265'
$ws.Range("B2").Value = 'Advanced Calculus'
$ws.Range("F2").Value = 'Dr. Meena Patel'
$ws.Range("A3").Value = 'CS265'
$ws.Range("B3").Value = 'System Software'
$ws.Range("F3").Value = 'Dr. Rajeev Malhotra'
$ws.Range("A4").Value = 'CS268'
$ws.Range("B4").Value = 'Algorithm Design'
$ws.Range("F4").Value = 'Dr. Ananya Das'
$ws.Range("A5").Value = 'CS269'
$ws.Range("B5").Value = 'Network Systems'
$ws.Range("F5").Value = 'Dr. Ananya Das'
